# [PV-94][WIP] Support for plans without sticky-ids or levels
# Rename the import header row so the reader no longer depends on
# "Unique Sticky ID" / "Start" / "Finish" columns:
#   A1: Unique Sticky ID -> Row ID
#   C1: Name             -> Task
#   E1: Start            -> Start Date
#   F1: Finish           -> End Date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03-t02-activity-added")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

$ws.Range("F1").Select()
